$d = $word.ActiveDocument

# Row 1 time slot: 08:30 - 08:45 -> 08:45 - 09:00
$d.Content.Find.Execute("08:30 - 08:45", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "08:45 - 09:00", 2)

# Row 2 time slot: 09:00 - 09:15 -> 09:15 - 09:30
$d.Content.Find.Execute("09:00 - 09:15", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "09:15 - 09:30", 2)

# Row 3 time slot: 10:00 - 10:15 -> 09:30 - 09:45
$d.Content.Find.Execute("10:00 - 10:15", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "09:30 - 09:45", 2)

# Swap the sellers in row 1 and row 3 (INTERLINK2AMERICAS <-> BOX BRAND).
# Route through a placeholder so the second replace doesn't re-match the
# text just written by the first one.
$d.Content.Find.Execute("INTERLINK2AMERICAS", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "__TMP_SELLER__", 2)
$d.Content.Find.Execute("BOX BRAND", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "INTERLINK2AMERICAS", 2)
$d.Content.Find.Execute("__TMP_SELLER__", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "BOX BRAND", 2)
